# Insert a new data row before row 123 (this shifts the existing rows
# 123..170 down to 124..171, growing the sheet's used range to A1:R171)
# and then populate the newly inserted row with the new weekly record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(123).Insert()

$ws.Cells.Item(123, 1).Value = 4
$ws.Cells.Item(123, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(123, 3).Value = "Los Lagos"
$ws.Cells.Item(123, 4).Value = 44553
$ws.Cells.Item(123, 5).Value = 10
$ws.Cells.Item(123, 6).Value = 100112032
$ws.Cells.Item(123, 7).Value = "Zapallo italiano"
$ws.Cells.Item(123, 8).Value = "Sin especificar"
$ws.Cells.Item(123, 9).Value = "Primera"
$ws.Cells.Item(123, 10).Value = 160
$ws.Cells.Item(123, 11).Value = 11000
$ws.Cells.Item(123, 12).Value = 12000
$ws.Cells.Item(123, 13).Value = 11500
$ws.Cells.Item(123, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(123, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(123, 16).Value = 230
$ws.Cells.Item(123, 17).Value = 50
$ws.Cells.Item(123, 18).Value = "Hortaliza"
